$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.8722
$ws.Range("E21").Value = 12.9958
$ws.Range("E23").Value = 14.2084
$ws.Range("E25").Value = 13.1712
